$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new log entry row (row 25) following the same pattern as the
# preceding rows: Date | "Internship" | Description.
$newRow = 25

# Copy formatting from the row above so the new row matches existing styling
$ws.Range("A24:C24").Copy() | Out-Null
$ws.Range("A25:C25").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item($newRow, 1).Value = 45260
$ws.Cells.Item($newRow, 2).Value = "Internship"
$ws.Cells.Item($newRow, 3).Value = "Contributed technical work by aiding in resolving inconsistencies flagged by the system for employee calls"

$ws.Range("C26").Select() | Out-Null
